$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column M (shifts M,N,O -> N,O,P)
$ws.Columns("M").Insert()

# Set header for new column M
$ws.Range("M1").Value = "Gas Supplier"

# Set data value for new column M, row 2
$ws.Range("M2").Value = "N/A"
